$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert the three new columns (rightmost first so earlier column letters stay put
# while we work on them): a new column after the old G, after the old E, and after
# the old C. Inserting this way reproduces the target layout:
#   A B C D(new) E F G(new) H I J(new)
$ws.Range("H1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()

# Row 2 headers for the three newly inserted columns.
$ws.Range("D2").Value = "BF01"
$ws.Range("G2").Value = "BF01"
$ws.Range("J2").Value = "BF01"

# New formula columns: each new column is the reciprocal of the BF10 column
# immediately to its left.
$ws.Range("D3").Formula = "=1/C3"
$ws.Range("D4:D5").Formula = "=1/C4"

$ws.Range("G3").Formula = "=1/F3"
$ws.Range("G4:G5").Formula = "=1/F4"

$ws.Range("J3").Formula = "=1/I3"
$ws.Range("J4:J5").Formula = "=1/I4"

# Selection moves as recorded in the saved view state.
$ws.Range("F15").Select()
